$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# FTF / MTTC record id values were regenerated -> update the two changed ids
$ws.Range("D2").Value = "a0Nq0000003PKUc"
$ws.Range("D3").Value = "a0Nq0000003PKUh"

# Update selected/active cell to D3 (was W8)
$ws.Range("D3").Select()
